$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.392.00"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "1.747.32"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.12"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4811"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2616"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "1.743.40"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.07"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06934"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6030"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.475"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.28"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9995"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "27.365.33"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007053"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "1.969.12"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.442"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.428"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.097"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.99"
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.27"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.841"
$ws.Range("E27").Value = "  +5.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.41"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.954"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07978"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.670"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04685"
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.597"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.012"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6178"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9241"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.546"
$ws.Range("E38").Value = "  +6.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.000"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9991"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.738"
$ws.Range("E41").Value = "  +5.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01495"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.88"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3836"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.876"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05362"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.863"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.86"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.245"
$ws.Range("E50").Value = "  +2.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.00"
$ws.Range("E51").Value = "  -0.48%  "
